$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44167
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 19000
$ws.Range("P2").Value = 18500
$ws.Range("Q2").Value = "`$/caja 13 kilos"
$ws.Range("S2").Value = 1423
$ws.Range("T2").Value = 13

# Row 3
$ws.Range("D3").Value = 44482
$ws.Range("L3").Value = "Primera"
$ws.Range("N3").Value = 25000
$ws.Range("O3").Value = 26000
$ws.Range("P3").Value = 25500
$ws.Range("S3").Value = 2125

# Row 4
$ws.Range("D4").Value = 44489
$ws.Range("N4").Value = 24000
$ws.Range("O4").Value = 25000
$ws.Range("P4").Value = 24500
$ws.Range("Q4").Value = "`$/caja 12 kilos"
$ws.Range("S4").Value = 2042

# Row 5
$ws.Range("D5").Value = 45125
$ws.Range("M5").Value = 160
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 14375
$ws.Range("Q5").Value = "`$/bandeja 10 kilos"
$ws.Range("S5").Value = 1438
$ws.Range("T5").Value = 10

# Row 6
$ws.Range("D6").Value = 45125
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 180
$ws.Range("N6").Value = 13000
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 13000
$ws.Range("Q6").Value = "`$/bandeja 10 kilos"
$ws.Range("R6").Value = "Región de Coquimbo"
$ws.Range("S6").Value = 1300
$ws.Range("T6").Value = 10

# Row 7
$ws.Range("D7").Value = 45126
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 160
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 14375
$ws.Range("Q7").Value = "`$/bandeja 10 kilos"
$ws.Range("S7").Value = 1438
$ws.Range("T7").Value = 10

# Row 8
$ws.Range("D8").Value = 45126
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 180
$ws.Range("N8").Value = 13000
$ws.Range("O8").Value = 13000
$ws.Range("P8").Value = 13000
$ws.Range("Q8").Value = "`$/bandeja 10 kilos"
$ws.Range("S8").Value = 1300
$ws.Range("T8").Value = 10

# Row 9
$ws.Range("D9").Value = 44811
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 29000
$ws.Range("O9").Value = 30000
$ws.Range("P9").Value = 29500
$ws.Range("S9").Value = 2458

# Row 10
$ws.Range("D10").Value = 44860
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 23000
$ws.Range("O10").Value = 24000
$ws.Range("P10").Value = 23500
$ws.Range("R10").Value = "Provincia de Limarí"
$ws.Range("S10").Value = 1958

# Row 11
$ws.Range("D11").Value = 44881
$ws.Range("N11").Value = 22000
$ws.Range("O11").Value = 23000
$ws.Range("P11").Value = 22500
$ws.Range("S11").Value = 1875

# Row 12
$ws.Range("D12").Value = 44545
$ws.Range("M12").Value = 200
$ws.Range("N12").Value = 23000
$ws.Range("O12").Value = 24000
$ws.Range("P12").Value = 23500
$ws.Range("Q12").Value = "`$/bandeja 12 kilos"
$ws.Range("S12").Value = 1958
$ws.Range("T12").Value = 12

# Row 13
$ws.Range("D13").Value = 45133
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 150
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 16000
$ws.Range("P13").Value = 15667
$ws.Range("S13").Value = 1567

# Row 14
$ws.Range("D14").Value = 44783
$ws.Range("L14").Value = "Tercera"
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 27000
$ws.Range("O14").Value = 28000
$ws.Range("P14").Value = 27500
$ws.Range("S14").Value = 2292

# Row 15
$ws.Range("D15").Value = 44496
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 23000
$ws.Range("O15").Value = 24000
$ws.Range("P15").Value = 23500
$ws.Range("S15").Value = 1958

# Row 16
$ws.Range("D16").Value = 44524
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = 23000
$ws.Range("O16").Value = 24000
$ws.Range("P16").Value = 23500
$ws.Range("S16").Value = 1958

# Row 18
$ws.Range("D18").Value = 45205
$ws.Range("M18").Value = 200
$ws.Range("N18").Value = 22000
$ws.Range("O18").Value = 23000
$ws.Range("P18").Value = 22500
$ws.Range("Q18").Value = "`$/bandeja 10 kilos"
$ws.Range("S18").Value = 2250

# Row 19
$ws.Range("D19").Value = 44839
$ws.Range("M19").Value = 160
$ws.Range("N19").Value = 26000
$ws.Range("O19").Value = 27000
$ws.Range("P19").Value = 26500
$ws.Range("Q19").Value = "`$/caja 12 kilos"
$ws.Range("S19").Value = 2208
$ws.Range("T19").Value = 12

# Row 20
$ws.Range("D20").Value = 44846
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 160
$ws.Range("N20").Value = 24000
$ws.Range("O20").Value = 25000
$ws.Range("P20").Value = 24500
$ws.Range("Q20").Value = "`$/caja 12 kilos"
$ws.Range("S20").Value = 2042
$ws.Range("T20").Value = 12

# Row 21
$ws.Range("D21").Value = 44846
$ws.Range("L21").Value = "Segunda"
$ws.Range("M21").Value = 100
$ws.Range("Q21").Value = "`$/caja 12 kilos"
$ws.Range("S21").Value = 1875
$ws.Range("T21").Value = 12

# Row 22
$ws.Range("D22").Value = 44441
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = 29000
$ws.Range("O22").Value = 30000
$ws.Range("P22").Value = 29500
$ws.Range("S22").Value = 2458

# Row 23
$ws.Range("D23").Value = 45147
$ws.Range("M23").Value = 270
$ws.Range("N23").Value = 17000
$ws.Range("O23").Value = 18000
$ws.Range("P23").Value = 17500
$ws.Range("Q23").Value = "`$/caja 10 kilos"
$ws.Range("S23").Value = 1750

# Row 25
$ws.Range("D25").Value = 44160
$ws.Range("L25").Value = "Segunda"
$ws.Range("M25").Value = 200
$ws.Range("N25").Value = 19000
$ws.Range("O25").Value = 20000
$ws.Range("P25").Value = 19500
$ws.Range("Q25").Value = "`$/caja 13 kilos"
$ws.Range("S25").Value = 1500
$ws.Range("T25").Value = 13

# Row 26
$ws.Range("D26").Value = 44468
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 200
$ws.Range("N26").Value = 29000
$ws.Range("O26").Value = 30000
$ws.Range("P26").Value = 29500
$ws.Range("S26").Value = 2950

# Row 27
$ws.Range("D27").Value = 44874
$ws.Range("L27").Value = "Segunda"
$ws.Range("M27").Value = 250
$ws.Range("N27").Value = 22000
$ws.Range("O27").Value = 23000
$ws.Range("P27").Value = 22500
$ws.Range("S27").Value = 1875
